# Apply cryptos list update (values refreshed by scheduled GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "44.326.80"
$ws.Range("E2").Value2 = "  +3.50%  "
$ws.Range("D3").Value2 = "2.271.84"
$ws.Range("E3").Value2 = "  +2.70%  "
$ws.Range("E4").Value2 = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "323.47"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  +2.13%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "105.32"
$c.Style = "Normal"
$ws.Range("E6").Value2 = "  +6.33%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "0.591"
$c.Style = "Normal"
$ws.Range("E7").Value2 = "  +0.42%  "
$ws.Range("E8").Value2 = "  +0.18%  "
$ws.Range("E9").Value2 = "  +1.97%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "38.81"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  +5.17%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = "0.0844"
$c.Style = "Normal"
$ws.Range("E11").Value2 = "  +2.12%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = "7.90"
$c.Style = "Normal"
$ws.Range("E12").Value2 = "  +2.73%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "0.885"
$c.Style = "Normal"
$ws.Range("E14").Value2 = "  +3.06%  "
$ws.Range("D15").Value2 = "2.618.93"
$ws.Range("E15").Value2 = "  +2.82%  "
$ws.Range("E16").Value2 = "  +2.68%  "
$ws.Range("D17").Value2 = "2.273.60"
$ws.Range("E17").Value2 = "  +2.78%  "
$ws.Range("D18").Value2 = "44.243.10"
$ws.Range("E18").Value2 = "  +3.35%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "13.86"
$c.Style = "Normal"
$ws.Range("E19").Value2 = "  -3.64%  "
$ws.Range("E20").Value2 = "  +4.47%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "6.53"
$c.Style = "Normal"
$ws.Range("E21").Value2 = "  +1.60%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = "66.44"
$c.Style = "Normal"
$ws.Range("E22").Value2 = "  +1.86%  "
$ws.Range("E23").Value2 = "  +2.20%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "240.82"
$c.Style = "Normal"
$ws.Range("E24").Value2 = "  +1.95%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "2.22"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  +4.33%  "
$ws.Range("E26").Value2 = "  -0.28%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "10.28"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  +3.03%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = "38.43"
$c.Style = "Normal"
$ws.Range("E28").Value2 = "  +12.28%  "
$ws.Range("E29").Value2 = "  -0.54%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = "6.50"
$c.Style = "Normal"
$ws.Range("E30").Value2 = "  +3.42%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "164.14"
$c.Style = "Normal"
$ws.Range("E31").Value2 = "  +6.54%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = "20.69"
$c.Style = "Normal"
$ws.Range("E32").Value2 = "  +0.89%  "
$ws.Range("E33").Value2 = "  -0.79%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "2.76"
$c.Style = "Normal"
$ws.Range("E34").Value2 = "  -0.54%  "
$ws.Range("E35").Value2 = "  +9.58%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "2.02"
$c.Style = "Normal"
$ws.Range("E36").Value2 = "  +4.98%  "
$ws.Range("E37").Value2 = "  +2.20%  "
$ws.Range("E38").Value2 = "  +0.36%  "
$ws.Range("E39").Value2 = "  +3.49%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "4.39"
$c.Style = "Normal"
$ws.Range("E40").Value2 = "  -0.47%  "
$ws.Range("B41").Value2 = "Celestia"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "15.59"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  +27.86%  "
$ws.Range("B42").Value2 = "VeChain"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "0.0329"
$c.Style = "Normal"
$ws.Range("E42").Value2 = "  +1.41%  "
$ws.Range("E43").Value2 = "  +0.26%  "
$ws.Range("D44").Value2 = "1.779.61"
$ws.Range("E44").Value2 = "  -1.78%  "
$ws.Range("E45").Value2 = "  +0.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = "86.29"
$c.Style = "Normal"
$ws.Range("E46").Value2 = "  -2.79%  "
$ws.Range("E47").Value2 = "  +1.16%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "60.72"
$c.Style = "Normal"
$ws.Range("E48").Value2 = "  -0.06%  "
$ws.Range("E49").Value2 = "  +10.35%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "75.49"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -0.03%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = "104.24"
$c.Style = "Normal"
$ws.Range("E51").Value2 = "  +1.19%  "
